# Fgf15-Klb.xlsx : update scripts with new tpm
#
# The underlying NATMI computation was re-run with updated TPM values.
# - The "Resolving-Mac" sending-cluster block (original rows 14-17) is
#   dropped entirely from the recomputed output.
# - Every numeric column (E:T) for the remaining rows (2-13) is
#   refreshed with the newly computed values. The text columns
#   (A sending cluster, B ligand symbol, C receptor symbol, D target
#   cluster) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Resolving-Mac" sender rows (14-17); the resized used
# range becomes A1:T13, matching the recomputed export.
$ws.Rows("14:17").Delete()

# Recomputed values for columns E (Ligand-expressing cells) through
# T (Edge total expression derived specificity), rows 2-13 in order.
$newValues = @(
    @(3, 1, 0.07580100000000001, 0.227403, 0.2364650501155793, 0.2364650501155794, 2, 0.6666666666666666, 0.4958466666666667, 1.48754, 0.1397696882980165, 0.1397696882980165, 0.03758567318000001, 0.3382710586200001, 0.03305064634802937, 0.03305064634802938),
    @(3, 1, 0.07580100000000001, 0.227403, 0.2364650501155793, 0.2364650501155794, 3, 1, 1.755888, 5.267664, 0.4949512317912007, 0.4949512317912008, 0.133098066288, 1.197882596592, 0.117038667830274, 0.117038667830274),
    @(3, 1, 0.07580100000000001, 0.227403, 0.2364650501155793, 0.2364650501155794, 3, 1, 1.175512, 3.526536, 0.3313543417264301, 0.3313543417264302, 0.08910498511200002, 0.8019448660080001, 0.07835372102235511, 0.07835372102235512),
    @(3, 1, 0.07580100000000001, 0.227403, 0.2364650501155793, 0.2364650501155794, 3, 1, 0.1203513333333333, 0.361054, 0.03392473818435271, 0.03392473818435272, 0.009122751417999999, 0.082104762762, 0.008022014914920871, 0.008022014914920876),
    @(1, 0.3333333333333333, 0.133329, 0.399987, 0.4159265533022002, 0.4159265533022002, 2, 0.6666666666666666, 0.4958466666666667, 1.48754, 0.1397696882980165, 0.1397696882980165, 0.06611074022000001, 0.59499666198, 0.05813392470991686, 0.05813392470991686),
    @(1, 0.3333333333333333, 0.133329, 0.399987, 0.4159265533022002, 0.4159265533022002, 3, 1, 1.755888, 5.267664, 0.4949512317912007, 0.4949512317912008, 0.234110791152, 2.106997120368, 0.2058633598915925, 0.2058633598915925),
    @(1, 0.3333333333333333, 0.133329, 0.399987, 0.4159265533022002, 0.4159265533022002, 3, 1, 1.175512, 3.526536, 0.3313543417264301, 0.3313543417264302, 0.156729839448, 1.410568555032, 0.1378190692759935, 0.1378190692759935),
    @(1, 0.3333333333333333, 0.133329, 0.399987, 0.4159265533022002, 0.4159265533022002, 3, 1, 0.1203513333333333, 0.361054, 0.03392473818435271, 0.03392473818435272, 0.016046322922, 0.144416906298, 0.01411019942469736, 0.01411019942469737),
    @(3, 1, 0.111429, 0.334287, 0.3476083965822204, 0.3476083965822204, 2, 0.6666666666666666, 0.4958466666666667, 1.48754, 0.1397696882980165, 0.1397696882980165, 0.05525169822000001, 0.49726528398, 0.04858511724007024, 0.04858511724007025),
    @(3, 1, 0.111429, 0.334287, 0.3476083965822204, 0.3476083965822204, 3, 1, 1.755888, 5.267664, 0.4949512317912007, 0.4949512317912008, 0.195656843952, 1.760911595568, 0.1720492040693342, 0.1720492040693342),
    @(3, 1, 0.111429, 0.334287, 0.3476083965822204, 0.3476083965822204, 3, 1, 1.175512, 3.526536, 0.3313543417264301, 0.3313543417264302, 0.130986126648, 1.178875139832, 0.1151815514280815, 0.1151815514280815),
    @(3, 1, 0.111429, 0.334287, 0.3476083965822204, 0.3476083965822204, 3, 1, 0.1203513333333333, 0.361054, 0.03392473818435271, 0.03392473818435272, 0.013410628722, 0.120695658498, 0.01179252384473447, 0.01179252384473448)
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $i + 2
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 5   # column E == 5
        $ws.Cells.Item($r, $col).Value2 = $rowVals[$j]
    }
}

"Fgf15-Klb sheet updated with recomputed TPM values"
